$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number, Coin (B), Link (C), Price (D), Volume 1h (E)
$rows = @(
    @{ Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='29.545.90'; E='  +2.33%  ' },
    @{ Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='1.988.31'; E='  +5.86%  ' },
    @{ Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.005'; E='  +0.23%  ' },
    @{ Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='329.45'; E='  +1.27%  ' },
    @{ Row=6; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.005'; E='  +0.32%  ' },
    @{ Row=7; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.4679'; E='  +1.73%  ' },
    @{ Row=8; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.3946'; E='  +1.69%  ' },
    @{ Row=9; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='46.58'; E='  +0.11%  ' },
    @{ Row=10; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.07970'; E='  +1.20%  ' },
    @{ Row=11; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='1.003'; E='  +1.91%  ' },
    @{ Row=12; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='22.81'; E='  +4.96%  ' },
    @{ Row=13; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='2.033.17'; E='  +8.55%  ' },
    @{ Row=14; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.263'; E='  +3.92%  ' },
    @{ Row=15; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.882'; E='  +3.97%  ' },
    @{ Row=16; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.07167'; E='  +2.93%  ' },
    @{ Row=17; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='88.89'; E='  +0.77%  ' },
    @{ Row=18; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.008'; E='  +0.45%  ' },
    @{ Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000009974'; E='  +0.01%  ' },
    @{ Row=20; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='17.37'; E='  +2.44%  ' },
    @{ Row=21; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.006'; E='  +0.35%  ' },
    @{ Row=22; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='29.625.73'; E='  +2.56%  ' },
    @{ Row=23; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='5.552'; E='  +5.78%  ' },
    @{ Row=24; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='11.29'; E='  +3.07%  ' },
    @{ Row=25; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='2.255.52'; E='  +7.86%  ' },
    @{ Row=26; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='2.123'; E='  +1.67%  ' },
    @{ Row=27; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='158.38'; E='  +1.61%  ' },
    @{ Row=28; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='19.67'; E='  +1.84%  ' },
    @{ Row=29; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='5.978'; E='  -0.45%  ' },
    @{ Row=30; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='120.45'; E='  +2.64%  ' },
    @{ Row=31; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.966'; E='  +1.98%  ' },
    @{ Row=32; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.09445'; E='  +1.08%  ' },
    @{ Row=33; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.8935'; E='  -1.14%  ' },
    @{ Row=34; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='5.289'; E='  +0.51%  ' },
    @{ Row=35; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='1.347'; E='  +2.21%  ' },
    @{ Row=36; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='3.191'; E='  -2.23%  ' },
    @{ Row=37; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.05841'; E='  +1.33%  ' },
    @{ Row=38; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.178'; E='  -0.78%  ' },
    @{ Row=39; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.02133'; E='  +3.08%  ' },
    @{ Row=40; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='7.926'; E='  +3.38%  ' },
    @{ Row=41; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.5764'; E='  +2.03%  ' },
    @{ Row=42; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1825'; E='  +3.38%  ' },
    @{ Row=43; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.000003116'; E='  +91.50%  ' },
    @{ Row=44; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='9.830'; E='  +1.79%  ' },
    @{ Row=45; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='12.16'; E='  +2.32%  ' },
    @{ Row=46; B='Decentraland'; C='https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; D='0.5385'; E='  +0.72%  ' },
    @{ Row=47; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='2.174'; E='  -3.32%  ' },
    @{ Row=48; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.644'; E='  +5.34%  ' },
    @{ Row=49; B='Cronos'; C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D='0.06967'; E='  -0.94%  ' },
    @{ Row=50; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='1.870'; E='  +1.43%  ' },
    @{ Row=51; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='114.77'; E='  +1.45%  ' }
)

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C

    # The Price column holds numbers formatted as plain text (e.g. "29.545.90",
    # "0.07970" with a significant trailing zero, etc). Whenever the text looks
    # like a single plain decimal number, Excel would otherwise silently coerce
    # it into a real number (dropping trailing zeros / changing precision), so we
    # force the cell to Text format first and restore the Normal style after.
    $dCell = $ws.Range("D" + $r.Row)
    $dVal = $r.D
    if ($dVal -match "^-?\d+(\.\d+)?$") {
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
        $dCell.Style = "Normal"
    } else {
        $dCell.Value = $dVal
    }

    $ws.Range("E" + $r.Row).Value = $r.E
}